# Slide 10, "Content Placeholder 2" shape: update the bullet list describing
# current issues and future directions.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shp = $s.Shapes.Item(2)

# 1) "Some client-server communication is failing." ->
#    "Some client-server communication is not yet working."
#    Delete the old paragraph and insert a fresh one with the new wording so
#    it stays a single clean run (rather than a diffed-up multi-run replace).
$tr = $shp.TextFrame.TextRange
$oldFailing = $tr.Paragraphs(4, 1)
$oldFailing.Delete()

$tr = $shp.TextFrame.TextRange
$afterUnimplemented = $tr.Paragraphs(3, 1)
$afterUnimplemented.InsertAfter("`rSome client-server communication is not yet working.") | Out-Null

# 2) Remove the "UNM App Contest" bullet that used to sit right under
#    "Future Directions:" (paragraph 7).
$tr = $shp.TextFrame.TextRange
$oldUnm = $tr.Paragraphs(7, 1)
$oldUnm.Delete()

# 3) Promote "Fix issues/bugs" (now paragraph 7) and
#    "Continue customizing the look and feel of our app" (now paragraph 8)
#    from lvl 3 up to lvl 2 (IndentLevel is 1-based, so 4 -> 3).
$tr = $shp.TextFrame.TextRange
$fixIssues = $tr.Paragraphs(7, 1)
$fixIssues.IndentLevel = 3

$tr = $shp.TextFrame.TextRange
$continueCustom = $tr.Paragraphs(8, 1)
$continueCustom.IndentLevel = 3

# 4) Remove "Have built-in safeguards for Spotify integration" (now paragraph 9).
$tr = $shp.TextFrame.TextRange
$oldSafeguards = $tr.Paragraphs(9, 1)
$oldSafeguards.Delete()

# 5) Re-add "UNM App Contest" at lvl 2, now right after "Continue customizing
#    the look and feel of our app" and before "Google Play Store".
$tr = $shp.TextFrame.TextRange
$continueCustom2 = $tr.Paragraphs(8, 1)
$continueCustom2.InsertAfter("`rUNM App Contest") | Out-Null
